$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Boulders")
$ws.Activate()

# New boulder added to the "Boulders" sheet: "Salter's Nick" at Shaftoe
$ws.Range("A10").Value = "Shaftoe"
$ws.Range("B10").Value = "Salter's Nick"
$ws.Range("C10").Value = 55.135770999999998
$ws.Range("D10").Value = -1.918345

# Reflect the new cursor position left after entering the row
$ws.Range("D11").Select()
